$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May 2020 to Jun 2020")
$ws.Name = "Jun 2020 to Jul 2020"
